$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player names (write order matters for shared-string table ordering)
$ws.Range("A3").Value = "Jerry"
$ws.Range("A2").Value = "Garrett"
$ws.Range("A4").Value = "Amy"
$ws.Range("A5").Value = "Carol"
$ws.Range("A6").Value = "Dina"
$ws.Range("A7").Value = "Sandra"

# Row 2 (Garrett) scores
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1

# Row 3 (Jerry) scores
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 1

# Row 4 (Amy) scores
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1

# Row 5 (Carol) scores
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1

# Row 6 (Dina) scores
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 1

# Row 7 (Sandra) scores
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 2

$ws.Range("B4").Select()
